$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pay Roll Report")

# Replace the Pay Type label cells (A7:A11) with plain numbers 1-5
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5

# Fill in row 14 "Total Expected" with values, black font
$ws.Range("A14").Value = "Total Expected"
$ws.Range("C14").Value = 168
$ws.Range("D14").Value = 168
$ws.Range("E14").Value = 168
$ws.Range("F14").Value = 168
$ws.Range("C14").Font.Color = 0
$ws.Range("D14").Font.Color = 0
$ws.Range("E14").Font.Color = 0
$ws.Range("F14").Font.Color = 0

# Fill in row 15 "Result" with values
$ws.Range("A15").Value = "Result"
$ws.Range("C15").Value = -143
$ws.Range("D15").Value = -150.5
$ws.Range("E15").Value = 32
$ws.Range("F15").Value = -92.5

# Negative results in red, positive result in black
$ws.Range("C15").Font.Color = 255
$ws.Range("D15").Font.Color = 255
$ws.Range("E15").Font.Color = 0
$ws.Range("F15").Font.Color = 255

# Update selection to match the final state
$ws.Range("G20").Select()
